$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Words" sheet — add the two new Wordle puzzles (id 227 "THOSE",
#    id 228 "MOIST") by extending the dated log two more rows.
# ---------------------------------------------------------------------------
$wsWords = $wb.Worksheets.Item("Words")

$wsWords.Range("A11").Formula = "=A10+1"
$wsWords.Range("A11").NumberFormat = "mm-dd-yy"
$wsWords.Range("B11").Formula = "=B10+1"
$wsWords.Range("C11").Value = "THOSE"
$wsWords.Range("C11").Font.Color = 13369344

$wsWords.Range("A12").Formula = "=A11+1"
$wsWords.Range("A12").NumberFormat = "mm-dd-yy"
$wsWords.Range("B12").Formula = "=B11+1"
$wsWords.Range("C12").Value = "MOIST"
$wsWords.Range("C12").Font.Color = 13369344

# ---------------------------------------------------------------------------
# 2) "Results" sheet — log everyone's attempts for wordle 227 (rows 44-49)
#    and wordle 228 (rows 50-55). Each player's row repeats the pattern
#    already used throughout the sheet: date/id/name mirror the row six
#    above (previous day for the same player), "Word" is looked up via the
#    existing XLOOKUP array formula, and "Number of Attempts" is the
#    reported score. Phil has not yet reported his 228 score.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Results")

$players = @(44, 45, 46, 47, 48, 49)
foreach ($r in $players) {
    $prev = $r - 6
    $ws.Range("A$r").Formula = "=A$prev+1"
    $ws.Range("A$r").NumberFormat = "mm-dd-yy"
    $ws.Range("B$r").Formula = "=B$prev+1"
    $ws.Range("C$r").Formula = "=C$prev"
    $ws.Range("C$r").NumberFormat = "mm-dd-yy"
    $ws.Range("D$r").Formula = "=IF(XLOOKUP(`$B$r,wordle_ids,wordles)=`"`",`"`",XLOOKUP(`$B$r,wordle_ids,wordles))"
}

$attempts227 = @{44 = 2; 45 = 2; 46 = 3; 47 = 4; 48 = 3; 49 = 3}
foreach ($r in $attempts227.Keys) {
    $ws.Range("E$r").Value = $attempts227[$r]
}

foreach ($r in @(50, 51, 52, 53, 54, 55)) {
    $prev = $r - 6
    $ws.Range("A$r").Formula = "=A$prev+1"
    $ws.Range("A$r").NumberFormat = "mm-dd-yy"
    $ws.Range("B$r").Formula = "=B$prev+1"
    $ws.Range("C$r").Formula = "=C$prev"
    $ws.Range("C$r").NumberFormat = "mm-dd-yy"
    $ws.Range("D$r").Formula = "=IF(XLOOKUP(`$B$r,wordle_ids,wordles)=`"`",`"`",XLOOKUP(`$B$r,wordle_ids,wordles))"
}

# Phil (row 50) has no reported score yet for wordle 228.
$ws.Range("E51").Value = 5
$ws.Range("E52").Value = 3
$ws.Range("E53").Value = 4
$ws.Range("E54").Value = 4
$ws.Range("E55").Value = 6

# ---------------------------------------------------------------------------
# 3) View state — move the frozen-pane scroll position / selections to
#    follow the newly-entered data, same as the author would after typing.
# ---------------------------------------------------------------------------
$wsWords.Range("C12").Select()

$ws.Activate()
$ws.Range("E51").Select()
